$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "해외대학이 대기업 취직에 유리한 이유 (4)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/foreign-degree-job-market-merits-4/#utm_source=rss&utm_medium=rss&utm_campaign=foreign-degree-job-market-merits-4"

$ws.Range("D28").Value = "지능 이야기 #1"
$ws.Range("E28").Value = "https://ropiens.tistory.com/116"

$ws.Range("D41").Value = "오픈 포맷 데이터 성능 향상 전략 – Part. 2: Parquet 프로세싱"
$ws.Range("E41").Value = "http://cloudinsight.net/data/%ec%98%a4%ed%94%88-%ed%8f%ac%eb%a7%b7-%eb%8d%b0%ec%9d%b4%ed%84%b0-%ec%84%b1%eb%8a%a5-%ed%96%a5%ec%83%81-%ec%a0%84%eb%9e%b5-part-2-parquet-%ed%94%84%eb%a1%9c%ec%84%b8%ec%8b%b1/"

$ws.Range("D46").Value = "[Bioinformatics] 2021년 05월, 한국바이오협회 유전체 데이터 활용 머신러닝, 딥러닝 분석 교육-변이 판독기 학습 등 재직자 과정 교육생 모집"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/393"
